# Update odds values in Sheet1 to match the 2024-10-28 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Ludogorets - CSKA 1948 Sofia)
$ws.Range("G2").Value = 1.33
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 7.5
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("W2").Value = 4.75
$ws.Range("X2").Value = 5
$ws.Range("Z2").Value = 7.5
$ws.Range("AC2").Value = 7.5
$ws.Range("AH2").Value = 19
$ws.Range("AJ2").Value = 34
$ws.Range("AM2").Value = 126
$ws.Range("AN2").Value = 3
$ws.Range("AV2").Value = 101

# Row 3 (Metta - RFS)
$ws.Range("O3").Value = 1.07

# Row 5 (Basaksehir - Eyupspor)
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 3.5
$ws.Range("L5").Value = 4
$ws.Range("U5").Value = 1.73
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 8
$ws.Range("X5").Value = 10
$ws.Range("AG5").Value = 201
$ws.Range("AJ5").Value = 13
$ws.Range("AL5").Value = 29
$ws.Range("AN5").Value = 4
$ws.Range("AS5").Value = 151
$ws.Range("BC5").Value = 251

# Row 6 (Galatasaray - Besiktas)
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 4.1
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 21
$ws.Range("AD6").Value = 8.5
$ws.Range("AE6").Value = 11
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 15
$ws.Range("AM6").Value = 26
$ws.Range("AN6").Value = 4.33
